$d = $word.ActiveDocument
$dash = [char]0x2013

# --- Step 1 ---------------------------------------------------------------
# The "mixed" bullet currently reads (as 3 runs after the "--" run):
#   "mixed HEAD~1"  +  " or "  +  "git reset HEAD~1 <dash><nbsp>"
# We want it to read:
#   "mixed HEAD~1 " (bold+italic) + "<dash>" (italic only) + "<nbsp>" (bold+italic)
#
# First, delete " or git reset HEAD~1" (note: NOT the trailing space before
# the dash) so the run that used to hold "git reset HEAD~1 <dash><nbsp>"
# shrinks down to " <dash><nbsp>" while keeping its original bold+italic
# (b/bCs/i/iCs) run formatting untouched. Because that trimmed remainder now
# carries formatting that is byte-identical to the preceding "mixed HEAD~1"
# run, the COM host's own run-coalescing merges them into a single run:
#   "mixed HEAD~1 <dash><nbsp>"   (bold + italic, one <w:r>)
# The "--" run further to the left is never touched, so it stays separate,
# exactly like the unedited hunk context in the diff shows.
$r = $d.Content
$found = $r.Find.Execute(" or git reset HEAD~1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the ' or git reset HEAD~1' text to remove."
}
$r.Text = ""

# --- Step 2 ---------------------------------------------------------------
# Split the dash back out into its own run with italic-only formatting
# (matching the sibling "soft HEAD~1" bullet's "<dash><nbsp>" run), leaving
# the trailing non-breaking space in its own bold+italic run.
# Scope the search to start where step 1 left off so we hit *this* bullet's
# dash and not the "soft HEAD~1 <dash>" one earlier in the document.
$r2 = $d.Range($r.Start, $d.Content.End)
$found2 = $r2.Find.Execute($dash, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the dash character to reformat."
}
$dashRange = $d.Range($r2.Start, $r2.Start + 1)
$dashRange.Bold = 0
